$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 333334660
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 333334660
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 333334660
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -333335010
$ws.Range("H98").Value = 1805.5
$ws.Range("I98").Value = 881.875
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 881.875
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = 616.125
$ws.Range("N98").Value = -8496
$ws.Range("H122").Value = 1805.5
$ws.Range("I122").Value = 881.875
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 2645.625
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -195.625
$ws.Range("N122").Value = -21400
$ws.Range("H138").Value = 4939
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4939
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 14817
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -25097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4170.909
$ws.Range("I2").Value = 2037.5
$ws.Range("J2").Value = 9860
$ws.Range("K2").Value = 2037.5
$ws.Range("L2").Value = 9860
$ws.Range("M2").Value = -1924.5
$ws.Range("N2").Value = -10086
$ws.Range("H7").Value = 29740
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 29740
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 29740
$ws.Range("N7").Value = -29968
$ws.Range("H45").Value = 5452.273
$ws.Range("I45").Value = 4008.8572
$ws.Range("J45").Value = 7978.25
$ws.Range("K45").Value = 4008.8572
$ws.Range("L45").Value = 7978.25
$ws.Range("M45").Value = -3631.8572
$ws.Range("N45").Value = -8732.25
$ws.Range("H52").Value = 15078.182
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 15078.182
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15078.182
$ws.Range("N52").Value = -15714.182
$ws.Range("H74").Value = 1022.11426
$ws.Range("I74").Value = 962.4666999999999
$ws.Range("J74").Value = 1380
$ws.Range("K74").Value = 962.4666999999999
$ws.Range("L74").Value = 1380
$ws.Range("M74").Value = -88.46669999999995
$ws.Range("N74").Value = -3128
$ws.Range("H77").Value = 1022.11426
$ws.Range("I77").Value = 962.4666999999999
$ws.Range("J77").Value = 1380
$ws.Range("K77").Value = 4812.3335
$ws.Range("L77").Value = 6900
$ws.Range("M77").Value = -444.3334999999997
$ws.Range("N77").Value = -15636
$ws.Range("H110").Value = 6379.5835
$ws.Range("I110").Value = 6680.9443
$ws.Range("J110").Value = 5475.5
$ws.Range("K110").Value = 6680.9443
$ws.Range("L110").Value = 5475.5
$ws.Range("M110").Value = -4635.9443
$ws.Range("N110").Value = -9565.5
$ws.Range("H116").Value = 4170.909
$ws.Range("I116").Value = 2037.5
$ws.Range("J116").Value = 9860
$ws.Range("K116").Value = 2037.5
$ws.Range("L116").Value = 9860
$ws.Range("M116").Value = 256.5
$ws.Range("N116").Value = -14448
$ws.Range("H122").Value = 2145.889
$ws.Range("I122").Value = 2164.125
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6492.375
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4042.375
$ws.Range("N122").Value = -10900
$ws.Range("H127").Value = 42972.145
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 42972.145
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 42972.145
$ws.Range("N127").Value = -52892.145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4170.909
$ws.Range("I3").Value = 2037.5
$ws.Range("J3").Value = 9860
$ws.Range("K3").Value = 2037.5
$ws.Range("L3").Value = 9860
$ws.Range("M3").Value = -1923.5
$ws.Range("N3").Value = -10088
$ws.Range("H94").Value = 463.6
$ws.Range("I94").Value = 447.52942
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 447.52942
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = 3.470579999999984
$ws.Range("N94").Value = -1912
$ws.Range("H107").Value = 2198.9614
$ws.Range("I107").Value = 1809.1052
$ws.Range("J107").Value = 3257.1428
$ws.Range("K107").Value = 1809.1052
$ws.Range("L107").Value = 3257.1428
$ws.Range("M107").Value = 110.8948
$ws.Range("N107").Value = -7097.1428
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 48500
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 48500
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 48500
$ws.Range("N125").Value = -58340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 24693
$ws.Range("I108").Value = 12000
$ws.Range("J108").Value = 27866.25
$ws.Range("K108").Value = 12000
$ws.Range("L108").Value = 27866.25
$ws.Range("M108").Value = -8160
$ws.Range("N108").Value = -35546.25
$ws.Range("H122").Value = 1245.65
$ws.Range("I122").Value = 1030.2354
$ws.Range("J122").Value = 2466.3333
$ws.Range("K122").Value = 3090.7062
$ws.Range("L122").Value = 7398.999899999999
$ws.Range("M122").Value = -640.7062000000001
$ws.Range("N122").Value = -12298.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 834195.5
$ws.Range("I5").Value = 860.6
$ws.Range("J5").Value = 2223087
$ws.Range("K5").Value = 2581.8
$ws.Range("L5").Value = 6669261
$ws.Range("M5").Value = -2469.8
$ws.Range("N5").Value = -6669485
$ws.Range("H135").Value = 834195.5
$ws.Range("I135").Value = 860.6
$ws.Range("J135").Value = 2223087
$ws.Range("K135").Value = 7745.400000000001
$ws.Range("L135").Value = 20007783
$ws.Range("M135").Value = -5210.400000000001
$ws.Range("N135").Value = -20012853

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1439.4286
$ws.Range("I97").Value = 1195.2
$ws.Range("J97").Value = 2050
$ws.Range("K97").Value = 1195.2
$ws.Range("L97").Value = 2050
$ws.Range("M97").Value = -699.2
$ws.Range("N97").Value = -3042
$ws.Range("H102").Value = 2158.7856
$ws.Range("I102").Value = 1874.909
$ws.Range("J102").Value = 3199.6667
$ws.Range("K102").Value = 1874.909
$ws.Range("L102").Value = 3199.6667
$ws.Range("M102").Value = -252.9090000000001
$ws.Range("N102").Value = -6443.6667
$ws.Range("H113").Value = 1236
$ws.Range("I113").Value = 1230.1428
$ws.Range("J113").Value = 1256.5
$ws.Range("K113").Value = 1230.1428
$ws.Range("L113").Value = 1256.5
$ws.Range("M113").Value = 939.8571999999999
$ws.Range("N113").Value = -5596.5
$ws.Range("H122").Value = 2630.05
$ws.Range("I122").Value = 3939.4
$ws.Range("J122").Value = 2193.6
$ws.Range("K122").Value = 11818.2
$ws.Range("L122").Value = 6580.799999999999
$ws.Range("M122").Value = -9368.200000000001
$ws.Range("N122").Value = -11480.8
$ws.Range("H132").Value = 2577.1724
$ws.Range("I132").Value = 2069.5
$ws.Range("J132").Value = 3705.3333
$ws.Range("K132").Value = 6208.5
$ws.Range("L132").Value = 11115.9999
$ws.Range("M132").Value = -3678.5
$ws.Range("N132").Value = -16175.9999
$ws.Range("H133").Value = 47851.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 47851.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 47851.8
$ws.Range("N133").Value = -57971.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14287914
$ws.Range("I7").Value = 33334934
$ws.Range("J7").Value = 2649.75
$ws.Range("K7").Value = 33334934
$ws.Range("L7").Value = 2649.75
$ws.Range("M7").Value = -33334822
$ws.Range("N7").Value = -2873.75
$ws.Range("H40").Value = 2614.2856
$ws.Range("I40").Value = 1575
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1575
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1439
$ws.Range("N40").Value = -4272
$ws.Range("H68").Value = 2005.8889
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 2726.5
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2726.5
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -4224.5
$ws.Range("H71").Value = 2005.8889
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 2726.5
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 13632.5
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -21120.5
$ws.Range("H126").Value = 14287914
$ws.Range("I126").Value = 33334934
$ws.Range("J126").Value = 2649.75
$ws.Range("K126").Value = 100004802
$ws.Range("L126").Value = 7949.25
$ws.Range("M126").Value = -100002332
$ws.Range("N126").Value = -12889.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3533.3333
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -4648
$ws.Range("H65").Value = 3533.3333
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -23240
$ws.Range("H107").Value = 1266.2858
$ws.Range("I107").Value = 845.1
$ws.Range("J107").Value = 1827.8667
$ws.Range("K107").Value = 2535.3
$ws.Range("L107").Value = 5483.6001
$ws.Range("M107").Value = -615.3000000000002
$ws.Range("N107").Value = -9323.6001
$ws.Range("H119").Value = 26400
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 26400
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 26400
$ws.Range("N119").Value = -36076
$ws.Range("H126").Value = 1612.9354
$ws.Range("I126").Value = 1565.6957
$ws.Range("J126").Value = 1748.75
$ws.Range("K126").Value = 4697.0871
$ws.Range("L126").Value = 5246.25
$ws.Range("M126").Value = -2227.0871
$ws.Range("N126").Value = -10186.25
